$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF2").Value = "maa://25251 (92.0), ***maa://21730 (22.86), ***maa://39501 (18.18), *maa://36675 (60.0)"
$ws.Range("H3").Value = "maa://21247 (98.46), *maa://22748 (60.0)"
$ws.Range("X3").Value = "maa://27396 (84.62), maa://27484 (96.15), maa://27480 (82.86)"
$ws.Range("D6").Value = "maa://42407 (94.74)"
$ws.Range("L7").Value = "maa://28624 (93.55), maa://24957 (97.67)"
$ws.Range("AE7").NumberFormat = "@"
$ws.Range("AE7").Value = "4"
$ws.Range("AF7").Value = "*maa://26191 (67.9), *maa://36671 (68.0), *maa://42530 (60.0), maa://45272 (100.0)"
$ws.Range("A8").Value = "更新日期：2025.01.05 13:20:09"
$ws.Range("D8").Value = "*maa://21476 (73.47), **maa://39431 (50.0), *maa://37551 (57.14)"
$ws.Range("P8").Value = "maa://32931 (83.96), *maa://21916 (61.29), maa://23252 (92.42), maa://37496 (96.43), **maa://22759 (45.45)"
$ws.Range("X8").Value = "maa://21411 (95.95)"
$ws.Range("AF9").Value = "maa://26206 (89.52), *maa://22865 (51.92)"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "6"
$ws.Range("D10").Value = "***maa://25695 (19.34), **maa://32237 (40.48), ***maa://34206 (21.74), ***maa://39951 (16.28), ***maa://39243 (28.57), maa://45271 (100.0)"
$ws.Range("L10").Value = "**maa://24395 (39.29)"
$ws.Range("X10").Value = "maa://22301 (97.66), maa://22726 (100.0)"
$ws.Range("T11").Value = "maa://22747 (92.81), maa://22501 (98.61)"
$ws.Range("X12").Value = "maa://22753 (91.46), *maa://21485 (77.04), maa://37962 (86.67)"
$ws.Range("AF12").Value = "*maa://28932 (78.2), *maa://20106 (63.96), *maa://22769 (64.29)"
$ws.Range("D13").Value = "maa://24999 (91.82), maa://36673 (92.65), maa://25001 (85.51)"
$ws.Range("P13").Value = "maa://22676 (91.96), *maa://22583 (74.24), *maa://22500 (57.78)"
$ws.Range("AF13").Value = "**maa://22737 (32.37), maa://39883 (90.91), *maa://39885 (56.0)"
$ws.Range("AF15").Value = "maa://21364 (81.03), *maa://22766 (70.27), *maa://36666 (79.07)"
$ws.Range("AB19").Value = "*maa://30709 (63.68), *maa://36668 (56.41)"
$ws.Range("L20").Value = "maa://41331 (84.91)"
$ws.Range("AB21").Value = "*maa://21443 (80.0), ***maa://23820 (29.82)"
$ws.Range("AF21").Value = "maa://22524 (94.63), *maa://22432 (77.42)"
$ws.Range("X22").Value = "maa://21282 (98.48), *maa://37649 (66.67)"
$ws.Range("L23").Value = "maa://39756 (94.38), maa://39875 (93.85)"
$ws.Range("P23").Value = "maa://30587 (92.06), *maa://29748 (75.78), ***maa://29785 (16.42), *maa://37566 (70.0)"
$ws.Range("X24").Value = "maa://29988 (86.46), maa://23504 (93.18), **maa://22892 (39.31), *maa://25141 (76.56), *maa://36663 (78.26), ***maa://22815 (23.08)"
$ws.Range("AF24").Value = "maa://22523 (85.71), maa://36672 (80.77), maa://29910 (92.86), **maa://21440 (34.55)"
$ws.Range("D25").Value = "maa://29753 (94.92)"
$ws.Range("AB26").Value = "maa://42235 (93.59)"
$ws.Range("T27").Value = "*maa://30624 (77.36)"
$ws.Range("T28").Value = "maa://23263 (95.0), *maa://29765 (61.54)"
$ws.Range("X28").Value = "maa://39929 (89.88), ***maa://39723 (14.29), maa://41749 (91.67)"
$ws.Range("P29").Value = "*maa://23168 (55.56), *maa://30050 (53.33)"
$ws.Range("AF29").Value = "*maa://24080 (68.87), *maa://42865 (78.95), ***maa://34960 (8.33)"
$ws.Range("AB30").Value = "maa://42979 (96.61), maa://45045 (100.0)"
$ws.Range("L31").Value = "maa://35926 (93.38), maa://36258 (83.84), *maa://43904 (72.73)"
$ws.Range("T32").Value = "maa://42859 (96.63), maa://41108 (87.76), maa://41238 (96.3)"
$ws.Range("P33").Value = "*maa://21956 (79.86), *maa://22730 (79.31)"
$ws.Range("T34").Value = "maa://24526 (93.63)"
$ws.Range("P37").Value = "maa://21280 (89.22), *maa://21239 (66.67)"
$ws.Range("AF38").Value = "maa://36697 (86.49)"
$ws.Range("H39").Value = "maa://25199 (84.82), maa://36670 (87.91), maa://30434 (89.39), ***maa://25036 (16.0), *maa://44165 (66.67), *maa://45059 (66.67)"
$ws.Range("H44").Value = "maa://29768 (97.86), maa://27728 (96.0)"
$ws.Range("H46").Value = "maa://35931 (92.68), maa://43901 (90.0)"
$ws.Range("H47").Value = "maa://27410 (96.21), maa://29661 (97.86), maa://28038 (84.62)"
$ws.Range("H59").Value = "maa://27746 (83.18), maa://31270 (94.96)"
$ws.Range("H62").Value = "maa://42981 (96.88), maa://43903 (100.0)"
